$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new product (CONCOR AMLO 5/5 MG 30 TABS) was added to the shortage
# list as the first item. The previously-first item (ETHOXA ...) moves
# down to become the second item. Insert a new row right after the
# existing item row (row 7) to make room, duplicating row 7's layout
# (merged cells, styles, row height) for the item that is being pushed
# down.
$ws.Rows("8:8").Insert()

# Clone formatting + merged-cell layout of the item row (row 7) into the
# newly inserted row 8, then overwrite the values below.
$ws.Range("A7:Q7").Copy($ws.Range("A8"))
$ws.Rows("8:8").RowHeight = 24.75

# Helper: a handful of cells hold text that merely *looks* numeric
# ("1", "210.00", ...) even though their column's number format is a
# numeric one. Writing such a string straight into .Value lets the
# engine auto-coerce it into a real number (losing the literal
# formatting / shared-string-text nature). Briefly switching the cell
# to a text format while assigning the value keeps it a text value;
# restoring the original number format afterwards snaps the cell back
# onto its original style.
function Set-TextValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# Row 7 becomes the new item: CONCOR AMLO 5/5 MG 30 TABS
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "CONCOR AMLO 5/5 MG 30 TABS"
$ws.Range("H7").Value = "1:0"
Set-TextValue $ws.Range("L7") "1"
Set-TextValue $ws.Range("N7") "210.00"
Set-TextValue $ws.Range("P7") "69.3000"
$ws.Range("Q7").Value = "0:1"

# Row 8 becomes the item that used to be first: ETHOXA 250MG/5ML SYRUP 120ML
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "ETHOXA 250MG/5ML SYRUP 120ML"
$ws.Range("H8").Value = "0:0"
Set-TextValue $ws.Range("L8") "1"
Set-TextValue $ws.Range("N8") "99.00"
Set-TextValue $ws.Range("P8") "198.0000"
$ws.Range("Q8").Value = "2:0"

# Row 9 (previously row 8) holds the running total of the "sale price"
# column; update it to include the new item's price too, and match the
# taller row height used in the new layout.
$ws.Range("P9").Value = 267.3
$ws.Rows("9:9").RowHeight = 25.5

# Row 10 (previously row 9) is the footer with the report's generation
# timestamp; the time portion changed from 12:31 AM to 9:31 AM.
$ws.Range("A10").Value = "Wednesday, 18 June, 2025 9:31 AM"

Write-Output "done"
